# Insert a new timestamp column right before the static "nom"/"url_produit"
# columns (which were at EW/EX and shift right to EX/EY), mirroring the
# scraper appending a fresh price-history snapshot column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert blank column at EW; existing EW (nom) -> EX, EX (url_produit) -> EY.
$ws.Columns("EW:EW").Insert()

# New header cell: the scrape timestamp for this snapshot.
$ws.Range("EW1").Value = "2026-02-03 22:19:18"

# Determine the last data row (column A holds the product reference).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# For every product row, carry the latest known price (column EV) forward
# into the freshly inserted EW column when a price is present; rows whose
# price history has already gone blank stay blank.
for ($r = 2; $r -le $lastRow; $r++) {
    $priceCell = $ws.Cells.Item($r, 152)   # EV
    $price = $priceCell.Value2
    if ($price -is [double]) {
        $ws.Cells.Item($r, 153).Value2 = $price   # EW
    }
}
